# Auto-generated market-data refresh for Pandaemonium_Profits workbook
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H-N) per sheet/row
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 287.86365
$ws.Range("I33").Value = 249.35294
$ws.Range("J33").Value = 418.8
$ws.Range("K33").Value = 249.35294
$ws.Range("L33").Value = 418.8
$ws.Range("M33").Value = -20.35293999999999
$ws.Range("N33").Value = -876.8

$ws.Range("H62").Value = 2243.4285
$ws.Range("J62").Value = 2100
$ws.Range("L62").Value = 2100
$ws.Range("N62").Value = -3348

$ws.Range("H65").Value = 2243.4285
$ws.Range("J65").Value = 2100
$ws.Range("L65").Value = 10500
$ws.Range("N65").Value = -16740

$ws.Range("I86").Value = 168248.33
$ws.Range("J86").Value = 3660
$ws.Range("K86").Value = 168248.33
$ws.Range("L86").Value = 3660
$ws.Range("M86").Value = -167125.33
$ws.Range("N86").Value = -5906

$ws.Range("I89").Value = 168248.33
$ws.Range("J89").Value = 3660
$ws.Range("K89").Value = 841241.6499999999
$ws.Range("L89").Value = 18300
$ws.Range("M89").Value = -835625.6499999999
$ws.Range("N89").Value = -29532

$ws.Range("H92").Value = 584.1429000000001
$ws.Range("J92").Value = 552.5
$ws.Range("L92").Value = 552.5
$ws.Range("N92").Value = -3048.5

$ws.Range("H106").Value = 3307.52
$ws.Range("I106").Value = 3452.2354
$ws.Range("K106").Value = 3452.2354
$ws.Range("M106").Value = -2821.2354

$ws.Range("H112").Value = 1503.7059
$ws.Range("I112").Value = 800
$ws.Range("J112").Value = 1597.5333
$ws.Range("K112").Value = 2400
$ws.Range("L112").Value = 4792.5999
$ws.Range("M112").Value = -1292
$ws.Range("N112").Value = -7008.5999

$ws.Range("H125").Value = 5962.579
$ws.Range("I125").Value = 350
$ws.Range("J125").Value = 6274.3887
$ws.Range("K125").Value = 3150
$ws.Range("L125").Value = 56469.49830000001
$ws.Range("M125").Value = -690
$ws.Range("N125").Value = -61389.49830000001

$ws.Range("H137").Value = 679966.7
$ws.Range("I137").Value = 2326.0344
$ws.Range("J137").Value = 2317598.2
$ws.Range("K137").Value = 6978.1032
$ws.Range("L137").Value = 6952794.600000001
$ws.Range("M137").Value = -4428.1032
$ws.Range("N137").Value = -6957894.600000001

$ws.Range("H138").Value = 5223.696
$ws.Range("I138").Value = 1310.8485
$ws.Range("J138").Value = 15156.308
$ws.Range("K138").Value = 3932.5455
$ws.Range("L138").Value = 45468.924
$ws.Range("M138").Value = 1207.4545
$ws.Range("N138").Value = -55748.924

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 13357
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 13357
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 13357
$ws.Range("M32").ClearContents()
$ws.Range("N32").Value = -13931

$ws.Range("H61").Value = 4657.553
$ws.Range("I61").Value = 3179.8
$ws.Range("J61").Value = 13101.857
$ws.Range("K61").Value = 3179.8
$ws.Range("L61").Value = 13101.857
$ws.Range("M61").Value = -2967.8
$ws.Range("N61").Value = -13525.857

$ws.Range("H136").Value = 4657.553
$ws.Range("I136").Value = 3179.8
$ws.Range("J136").Value = 13101.857
$ws.Range("K136").Value = 9539.400000000001
$ws.Range("L136").Value = 39305.571
$ws.Range("M136").Value = -6989.400000000001
$ws.Range("N136").Value = -44405.571

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 748.5
$ws.Range("J99").Value = 497
$ws.Range("L99").Value = 497
$ws.Range("N99").Value = -3493

$ws.Range("H134").Value = 2232.3333
$ws.Range("I134").Value = 2102.0417
$ws.Range("K134").Value = 6306.125100000001
$ws.Range("M134").Value = -3771.125100000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 179.4
$ws.Range("J22").Value = 78
$ws.Range("L22").Value = 78
$ws.Range("N22").Value = -778

$ws.Range("H31").Value = 4346.6597
$ws.Range("I31").Value = 5225.846
$ws.Range("J31").Value = 3258.1428
$ws.Range("K31").Value = 5225.846
$ws.Range("L31").Value = 3258.1428
$ws.Range("M31").Value = -4930.846
$ws.Range("N31").Value = -3848.1428

$ws.Range("H34").Value = 4346.6597
$ws.Range("I34").Value = 5225.846
$ws.Range("J34").Value = 3258.1428
$ws.Range("K34").Value = 5225.846
$ws.Range("L34").Value = 3258.1428
$ws.Range("M34").Value = -5023.846
$ws.Range("N34").Value = -3662.1428

$ws.Range("H132").Value = 2023.1224
$ws.Range("I132").Value = 1613.4242
$ws.Range("J132").Value = 2868.125
$ws.Range("K132").Value = 4840.2726
$ws.Range("L132").Value = 8604.375
$ws.Range("M132").Value = -2310.2726
$ws.Range("N132").Value = -13664.375

$ws.Range("H134").Value = 3079.149
$ws.Range("I134").Value = 1718.2963
$ws.Range("K134").Value = 5154.8889
$ws.Range("M134").Value = -2619.8889

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 81502.414
$ws.Range("I131").Value = 2403.75
$ws.Range("J131").Value = 239699.75
$ws.Range("K131").Value = 7211.25
$ws.Range("L131").Value = 719099.25
$ws.Range("M131").Value = -2171.25
$ws.Range("N131").Value = -729179.25

$ws.Range("H132").Value = 1263.9111
$ws.Range("J132").Value = 1074.091
$ws.Range("L132").Value = 9666.819
$ws.Range("N132").Value = -14726.819

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2649.1538
$ws.Range("I132").Value = 2593.7144
$ws.Range("K132").Value = 7781.1432
$ws.Range("M132").Value = -5251.1432

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 746.5599999999999
$ws.Range("I22").Value = 601.5
$ws.Range("J22").Value = 880.46155
$ws.Range("K22").Value = 601.5
$ws.Range("L22").Value = 880.46155
$ws.Range("M22").Value = -306.5
$ws.Range("N22").Value = -1470.46155

$ws.Range("H27").Value = 746.5599999999999
$ws.Range("I27").Value = 601.5
$ws.Range("J27").Value = 880.46155
$ws.Range("K27").Value = 601.5
$ws.Range("L27").Value = 880.46155
$ws.Range("M27").Value = -494.5
$ws.Range("N27").Value = -1094.46155

$ws.Range("H46").Value = 1216.5
$ws.Range("I46").Value = 3000
$ws.Range("J46").Value = 859.8
$ws.Range("K46").Value = 3000
$ws.Range("L46").Value = 859.8
$ws.Range("M46").Value = -2812
$ws.Range("N46").Value = -1235.8

$ws.Range("H95").Value = 39866.668
$ws.Range("J95").Value = 39866.668
$ws.Range("L95").Value = 39866.668
$ws.Range("N95").Value = -45358.668

$ws.Range("H132").Value = 7274.724
$ws.Range("I132").Value = 9788.375
$ws.Range("J132").Value = 4181
$ws.Range("K132").Value = 29365.125
$ws.Range("L132").Value = 12543
$ws.Range("M132").Value = -26835.125
$ws.Range("N132").Value = -17603

$ws.Range("H136").Value = 4500.041
$ws.Range("I136").Value = 2418.3447
$ws.Range("J136").Value = 7518.5
$ws.Range("K136").Value = 7255.034100000001
$ws.Range("L136").Value = 22555.5
$ws.Range("M136").Value = -4705.034100000001
$ws.Range("N136").Value = -27655.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H137").Value = 44897.5
$ws.Range("J137").Value = 59795
$ws.Range("L137").Value = 59795
$ws.Range("N137").Value = -69995
